$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The TPM (transcripts-per-million) input used by the NATMI scripts changed,
# which updates the ligand average expression for the "ECs" sending cluster
# and the receptor average expression for the "ECs" target cluster. All the
# dependent specificity / edge-weight metrics in the table were recomputed
# from those two underlying values and are written out explicitly below.

# Row 2
$ws.Range("G2").Value = 35.906979
$ws.Range("H2").Value = 107.720937
$ws.Range("I2").Value = 0.6107087147789413
$ws.Range("J2").Value = 0.6107087147789412
$ws.Range("M2").Value = 201.098592
$ws.Range("N2").Value = 603.295776
$ws.Range("O2").Value = 0.7918622805845071
$ws.Range("P2").Value = 0.791862280584507
$ws.Range("Q2").Value = 7220.842919873568
$ws.Range("R2").Value = 64987.58627886211
$ws.Range("S2").Value = 0.4835971956576857
$ws.Range("T2").Value = 0.4835971956576856

# Row 3
$ws.Range("G3").Value = 35.906979
$ws.Range("H3").Value = 107.720937
$ws.Range("I3").Value = 0.6107087147789413
$ws.Range("J3").Value = 0.6107087147789412
$ws.Range("O3").Value = 0.1414593902976603
$ws.Range("P3").Value = 0.1414593902976603
$ws.Range("Q3").Value = 1289.941523829765
$ws.Range("R3").Value = 11609.47371446788
$ws.Range("S3").Value = 0.08639048244209674
$ws.Range("T3").Value = 0.08639048244209672

# Row 4
$ws.Range("G4").Value = 35.906979
$ws.Range("H4").Value = 107.720937
$ws.Range("I4").Value = 0.6107087147789413
$ws.Range("J4").Value = 0.6107087147789412
$ws.Range("O4").Value = 0.0666783291178327
$ws.Range("P4").Value = 0.06667832911783268
$ws.Range("Q4").Value = 608.02711850867
$ws.Range("R4").Value = 5472.24406657803
$ws.Range("S4").Value = 0.04072103667915886
$ws.Range("T4").Value = 0.04072103667915885

# Row 5
$ws.Range("I5").Value = 0.2899643113254147
$ws.Range("J5").Value = 0.2899643113254147
$ws.Range("M5").Value = 201.098592
$ws.Range("N5").Value = 603.295776
$ws.Range("O5").Value = 0.7918622805845071
$ws.Range("P5").Value = 0.791862280584507
$ws.Range("Q5").Value = 3428.454013805953
$ws.Range("R5").Value = 30856.08612425358
$ws.Range("S5").Value = 0.2296118008542589
$ws.Range("T5").Value = 0.2296118008542589

# Row 6
$ws.Range("I6").Value = 0.2899643113254147
$ws.Range("J6").Value = 0.2899643113254147
$ws.Range("O6").Value = 0.1414593902976603
$ws.Range("P6").Value = 0.1414593902976603
$ws.Range("S6").Value = 0.04101817468817411
$ws.Range("T6").Value = 0.04101817468817411

# Row 7
$ws.Range("I7").Value = 0.2899643113254147
$ws.Range("J7").Value = 0.2899643113254147
$ws.Range("O7").Value = 0.0666783291178327
$ws.Range("P7").Value = 0.06667832911783268
$ws.Range("S7").Value = 0.0193343357829817
$ws.Range("T7").Value = 0.0193343357829817

# Row 8
$ws.Range("I8").Value = 0.09932697389564409
$ws.Range("J8").Value = 0.09932697389564407
$ws.Range("M8").Value = 201.098592
$ws.Range("N8").Value = 603.295776
$ws.Range("O8").Value = 0.7918622805845071
$ws.Range("P8").Value = 0.791862280584507
$ws.Range("Q8").Value = 1174.413364096896
$ws.Range("R8").Value = 10569.72027687207
$ws.Range("S8").Value = 0.07865328407256253
$ws.Range("T8").Value = 0.07865328407256252

# Row 9
$ws.Range("I9").Value = 0.09932697389564409
$ws.Range("J9").Value = 0.09932697389564407
$ws.Range("O9").Value = 0.1414593902976603
$ws.Range("P9").Value = 0.1414593902976603
$ws.Range("S9").Value = 0.01405073316738943
$ws.Range("T9").Value = 0.01405073316738943

# Row 10
$ws.Range("I10").Value = 0.09932697389564409
$ws.Range("J10").Value = 0.09932697389564407
$ws.Range("O10").Value = 0.0666783291178327
$ws.Range("P10").Value = 0.06667832911783268
$ws.Range("S10").Value = 0.006622956655692134
$ws.Range("T10").Value = 0.006622956655692131

